$d = $word.ActiveDocument

# NOTE: Find.Execute's built-in "replace" (ReplaceWith) path runs typed
# text through smart-quote autocorrect, which mangles the straight
# apostrophes ( ' ) in a couple of the French strings below into curly
# ones. To keep the exact characters from the target translation, we
# instead use Find purely to *locate* each run of text, then assign the
# new string straight onto the matched Range (no autocorrect involved).

# The subtitle "[Music]" text appears twice in this document (the first
# and the very last cue). Only the first one is translated per the diff,
# so scope the Find to that specific paragraph's range rather than
# searching/replacing across the whole document.
$p6 = $d.Paragraphs.Item(6)
$r1 = $p6.Range
$r1.Find.Execute("[Music]", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r1.Text = "[Musique]"

# The remaining four subtitle lines are unique strings in the document, so
# a whole-document Find is safe for each of them.
$r2 = $d.Content
$r2.Find.Execute("four bright mathematicians are taken into", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r2.Text = "Quatre mathématiciens brillants sont placés en"

$r3 = $d.Content
$r3.Find.Execute("custody and put in jail because they tried", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r3.Text = "garde à vue et mis en prison parce qu'ils ont essayé"

$r4 = $d.Content
$r4.Find.Execute("to convince an old lady that the Goedel's", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r4.Text = "de convaincre une vieille dame que les"

$r5 = $d.Content
$r5.Find.Execute("incompleteness theorems are true. Every", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r5.Text = "théorèmes d'incomplétude de Goedel sont vrais. Chaque"
